# Remove the "STOCK"/"Capacity" data column from the PP worksheet's
# TOPOLOGY table (F5:F9) by deleting those cells and shifting G:K left.
# This breaks/re-targets the DMD!C18 formula that referenced PP!F8 (-> #REF!)
# and leaves the now-unused styles/borders/shared strings to be cleaned up
# on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PP")
$ws.Activate()

$ws.Range("F5:F9").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)

$ws.Application.Goto($ws.Range("F14"), $true)
